$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 6.400776119402984
$ws.Range("B2").Value = 2.903999999999999
$ws.Range("C2").Value = 10.852
$ws.Range("A3").Value = 0.08827860696517399
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0.4079999999999999
$ws.Range("A4").Value = 0.1794029850746263
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0.6280000000000006
$ws.Range("A5").Value = 0.1286766169154227
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0.4839999999999998
$ws.Range("A6").Value = 0.9091542288557202
$ws.Range("B6").Value = 0.188
$ws.Range("C6").Value = 1.859999999999999
$ws.Range("A7").Value = 4.390766169154213
$ws.Range("B7").Value = 1.783999999999999
$ws.Range("C7").Value = 7.904000000000003
$ws.Range("A8").Value = 11.07599999999999
$ws.Range("B8").Value = 6.036
$ws.Range("C8").Value = 16.41199999999999
$ws.Range("A9").Value = 0.7004577114427848
$ws.Range("B9").Value = 0.028
$ws.Range("C9").Value = 1.639999999999999
$ws.Range("A10").Value = 11.99291542288556
$ws.Range("B10").Value = 6.628000000000002
$ws.Range("C10").Value = 18.132
$ws.Range("A11").Value = 10.26722388059702
$ws.Range("B11").Value = 5.312000000000003
$ws.Range("C11").Value = 16.04799999999999
$ws.Range("A12").Value = 3.819422885572128
$ws.Range("B12").Value = 1.403999999999999
$ws.Range("C12").Value = 6.951999999999999
$ws.Range("A13").Value = 14.718407960199
$ws.Range("B13").Value = 8.084000000000001
$ws.Range("C13").Value = 22.308
$ws.Range("A14").Value = 12.53974129353232
$ws.Range("B14").Value = 6.896000000000002
$ws.Range("C14").Value = 18.628
$ws.Range("A15").Value = 0.1893532338308455
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = 0.6280000000000006
$ws.Range("A16").Value = 0.9099104477611928
$ws.Range("B16").Value = 0.18
$ws.Range("C16").Value = 1.879999999999999
$ws.Range("A17").Value = 3.759422885572131
$ws.Range("B17").Value = 1.451999999999999
$ws.Range("C17").Value = 6.764000000000002
$ws.Range("A18").Value = 0.1777711442786064
$ws.Range("B18").Value = 0
$ws.Range("C18").Value = 0.6520000000000006
$ws.Range("A19").Value = 12.6566567164179
$ws.Range("B19").Value = 6.836000000000004
$ws.Range("C19").Value = 19.144
$ws.Range("A20").Value = 1.251721393034825
$ws.Range("B20").Value = 0.308
$ws.Range("C20").Value = 2.571999999999999
$ws.Range("A21").Value = 0.3731343283582084
$ws.Range("B21").Value = 0
$ws.Range("C21").Value = 1.008
$ws.Range("A22").Value = 13.31661691542288
$ws.Range("B22").Value = 7.012
$ws.Range("C22").Value = 19.71199999999999
$ws.Range("A23").Value = 0.181054726368159
$ws.Range("B23").Value = 0
$ws.Range("C23").Value = 0.6400000000000005
$ws.Range("A24").Value = 0.9055522388059688
$ws.Range("B24").Value = 0.1640000000000001
$ws.Range("C24").Value = 1.879999999999999
$ws.Range("A25").Value = 0.5925572139303478
$ws.Range("B25").Value = 0
$ws.Range("C25").Value = 1.356
$ws.Range("A26").Value = 0.3521990049751241
$ws.Range("B26").Value = 0
$ws.Range("C26").Value = 1.008
$ws.Range("A27").Value = 16.89787064676617
$ws.Range("B27").Value = 9.383999999999999
$ws.Range("C27").Value = 24.96799999999999
$ws.Range("A28").Value = 5.4809751243781
$ws.Range("B28").Value = 2.212000000000002
$ws.Range("C28").Value = 9.444000000000004
$ws.Range("A29").Value = 0.226189054726367
$ws.Range("B29").Value = 0
$ws.Range("C29").Value = 0.78
$ws.Range("A30").Value = 3.19383084577113
$ws.Range("B30").Value = 1.02
$ws.Range("C30").Value = 6.127999999999999
$ws.Range("A31").Value = 11.40949253731342
$ws.Range("B31").Value = 5.840000000000004
$ws.Range("C31").Value = 17.792
$ws.Range("A32").Value = 11.2499502487562
$ws.Range("B32").Value = 6.083999999999998
$ws.Range("C32").Value = 17.24400000000001
$ws.Range("A33").Value = 12.7369751243781
$ws.Range("B33").Value = 7.044
$ws.Range("C33").Value = 19.01199999999999
$ws.Range("A34").Value = 8.952995024875607
$ws.Range("B34").Value = 4.391999999999999
$ws.Range("C34").Value = 14.628
$ws.Range("A35").Value = 1.092796019900496
$ws.Range("B35").Value = 0.108
$ws.Range("C35").Value = 2.404000000000001
$ws.Range("A36").Value = 4.5035422885572
$ws.Range("B36").Value = 1.783999999999999
$ws.Range("C36").Value = 8.02
$ws.Range("A37").Value = 0.442368159203979
$ws.Range("B37").Value = 0
$ws.Range("C37").Value = 1.204000000000001
$ws.Range("A38").Value = 14.34758208955224
$ws.Range("B38").Value = 8.307999999999995
$ws.Range("C38").Value = 21.24799999999999
$ws.Range("A39").Value = 4.076238805970144
$ws.Range("B39").Value = 1.608
$ws.Range("C39").Value = 7.688000000000006
$ws.Range("A40").Value = 0.2496318407960193
$ws.Range("B40").Value = 0
$ws.Range("C40").Value = 0.7999999999999998
$ws.Range("A41").Value = 3.750985074626853
$ws.Range("B41").Value = 1.484
$ws.Range("C41").Value = 6.824000000000004
$ws.Range("A42").Value = 11.27395024875621
$ws.Range("B42").Value = 6.071999999999998
$ws.Range("C42").Value = 17.41999999999999
$ws.Range("A43").Value = 3.19383084577113
$ws.Range("B43").Value = 1.02
$ws.Range("C43").Value = 6.127999999999999
$ws.Range("A44").Value = 6.849253731343281
$ws.Range("B44").Value = 3.424000000000003
$ws.Range("C44").Value = 11.192
$ws.Range("A45").Value = 13.34571144278606
$ws.Range("B45").Value = 7.644
$ws.Range("C45").Value = 20.19599999999999
$ws.Range("A46").Value = 10.11418905472636
$ws.Range("B46").Value = 5.595999999999999
$ws.Range("C46").Value = 15.612
$ws.Range("A47").Value = 14.09596019900495
$ws.Range("B47").Value = 8.1
$ws.Range("C47").Value = 21.45200000000001
$ws.Range("A48").Value = 3.873572139303471
$ws.Range("B48").Value = 1.448
$ws.Range("C48").Value = 7.024000000000001
$ws.Range("A49").Value = 0.2106268656716404
$ws.Range("B49").Value = 0
$ws.Range("C49").Value = 0.6320000000000005
$ws.Range("A50").Value = 0.4528159203980093
$ws.Range("B50").Value = 0
$ws.Range("C50").Value = 1.272
$ws.Range("A51").Value = 0.08382089552238794
$ws.Range("B51").Value = 0
$ws.Range("C51").Value = 0.3999999999999999
$ws.Range("A52").Value = 9.014567164179105
$ws.Range("B52").Value = 4.399999999999999
$ws.Range("C52").Value = 14.036
$ws.Range("A53").Value = 4.190228855721388
$ws.Range("B53").Value = 1.707999999999999
$ws.Range("C53").Value = 7.728000000000006
$ws.Range("A54").Value = 0.1725572139303481
$ws.Range("B54").Value = 0
$ws.Range("C54").Value = 0.6080000000000005
$ws.Range("A55").Value = 0.3579104477611936
$ws.Range("B55").Value = 0
$ws.Range("C55").Value = 1.008
$ws.Range("A56").Value = 11.72630845771144
$ws.Range("B56").Value = 5.988
$ws.Range("C56").Value = 17.776
$ws.Range("A57").Value = 11.25114427860695
$ws.Range("B57").Value = 6.095999999999997
$ws.Range("C57").Value = 17.24000000000001
$ws.Range("A58").Value = 0.4343084577114413
$ws.Range("B58").Value = 0
$ws.Range("C58").Value = 1.244
$ws.Range("A59").Value = 8.81653731343283
$ws.Range("B59").Value = 4.572000000000002
$ws.Range("C59").Value = 14.096
$ws.Range("A60").Value = 0.1964975124378106
$ws.Range("B60").Value = 0
$ws.Range("C60").Value = 0.6640000000000005
$ws.Range("A61").Value = 13.76859701492538
$ws.Range("B61").Value = 7.367999999999995
$ws.Range("C61").Value = 20.348
$ws.Range("A62").Value = 10.27932338308458
$ws.Range("B62").Value = 5.336000000000001
$ws.Range("C62").Value = 16.07999999999999
$ws.Range("A63").Value = 15.54680597014925
$ws.Range("B63").Value = 8.684000000000003
$ws.Range("C63").Value = 22.912
$ws.Range("A64").Value = 0.1992835820895518
$ws.Range("B64").Value = 0
$ws.Range("C64").Value = 0.6640000000000005
$ws.Range("A65").Value = 10.19590049751243
$ws.Range("B65").Value = 5.552000000000004
$ws.Range("C65").Value = 15.896
$ws.Range("A66").Value = 8.214467661691527
$ws.Range("B66").Value = 4.184000000000003
$ws.Range("C66").Value = 13.496
$ws.Range("A67").Value = 7.00784079601989
$ws.Range("B67").Value = 3.255999999999999
$ws.Range("C67").Value = 11.84
$ws.Range("A68").Value = 0.317890547263681
$ws.Range("B68").Value = 0
$ws.Range("C68").Value = 0.96
$ws.Range("A69").Value = 14.52039800995023
$ws.Range("B69").Value = 8.179999999999993
$ws.Range("C69").Value = 22.27599999999999
$ws.Range("A70").Value = 0.2562786069651734
$ws.Range("B70").Value = 0
$ws.Range("C70").Value = 0.7959999999999998
$ws.Range("A71").Value = 1.444855721393035
$ws.Range("B71").Value = 0.3999999999999999
$ws.Range("C71").Value = 3.159999999999999
$ws.Range("A72").Value = 2.273094527363169
$ws.Range("B72").Value = 0.7959999999999998
$ws.Range("C72").Value = 4.551999999999998
